$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Mark every still-unmarked image-carrying run as "do not spell-check"
#    (<w:noProof/>) -- this affects one run using Bold/Bold Bi styling and
#    three runs using the Arial body-text style (the two logo runs at the
#    top of the document already carry <w:noProof/> and are left alone).
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $xml = $para.Range.WordOpenXML
    if ($xml -match "<w:drawing>" -and -not ($xml -match "w:noProof")) {
        $para.Range.NoProofing = -1
    }
}

# ---------------------------------------------------------------------------
# 2. Locate the "English version:" placeholder paragraph together with the
#    stray empty paragraph right before it, merge the two (deleting the
#    empty paragraph's mark joins it into the following paragraph) and
#    replace the placeholder text with the final Spanish explanation of why
#    Google Drive (rather than the GitHub repo) was used for delivery.
# ---------------------------------------------------------------------------
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd([char]13)
    if ($text -eq "English version:") {
        if ($i -gt 1) {
            $prev = $d.Paragraphs.Item($i - 1)
            $prevText = $prev.Range.Text.TrimEnd([char]13)
            if ($prevText -eq "") {
                $emptyRange = $d.Range($prev.Range.Start, $prev.Range.End)
                $emptyRange.Delete()
            }
        }
        break
    }
}

$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $text = $para.Range.Text.TrimEnd([char]13)
    if ($text -eq "English version:") {
        $r1 = $para.Range
        $r1.Find.Execute("English ", $true, $false, $false, $false, $false, $true, 1, $false, `
            "Por ello, pese a que todo nuestro trabajo ", 2) | Out-Null

        $r2 = $para.Range
        $r2.Find.Execute("version", $true, $false, $false, $false, $false, $true, 1, $false, `
            "esta", 2) | Out-Null

        $r3 = $para.Range
        $r3.Find.Execute(":", $true, $false, $false, $false, $false, $true, 1, $false, `
            " realizado en github, decidimos usar Google Drive para la entrega final con el fin de resolver problemas de compatibilidad entre procesadores AMD Ryzen e Intel. ", 2) | Out-Null

        break
    }
}

Write-Output "done"
